# GenX FX Trading Dashboard - signal refresh
# Updates Active Signals, Summary Dashboard, and Signal History sheets
# to reflect the latest batch of trading signals.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: "Active Signals"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Active Signals")

# Every active signal is now a SELL, so make sure row 3-6's Signal cell
# (column C) carries the same red "SELL" fill/border as row 2, instead
# of the green "BUY" formatting some of them used before.
$sellFormat = $ws1.Range("C2")
$sellFormat.Copy()
$ws1.Range("C3:C6").PasteSpecial(-4122)

# Give the brand-new row 6 the same row-wide formatting (borders/fill)
# as the existing data rows before filling in its values.
$rowFormat = $ws1.Range("A2:J2")
$rowFormat.Copy()
$ws1.Range("A6:J6").PasteSpecial(-4122)

$activeSignals = @(
    @("2025-07-28 16:29", "NZDUSD", "SELL", 0.59027, 0.59471, 0.5813,    0.05, "68.0%", 2.02, "Active"),
    @("2025-07-28 16:03", "NZDUSD", "SELL", 0.59092, 0.59559, 0.58222,   0.01, "84.0%", 1.86, "Active"),
    @("2025-07-28 15:53", "USDJPY", "SELL", 148.7591, 148.97191, 148.00594, 0.08, "66.0%", 3.54, "Active"),
    @("2025-07-28 16:34", "NZDUSD", "SELL", 0.59032, 0.5938,  0.58405,   0.02, "87.0%", 1.8,  "Active"),
    @("2025-07-28 16:07", "USDCHF", "SELL", 0.88256, 0.88551, 0.87266,   0.09, "82.0%", 3.36, "Active")
)

# The Confidence column holds a literal text percentage (e.g. "68.0%"),
# not a numeric percentage, so force the column to Text format first -
# otherwise Excel auto-converts a "NN.N%" entry into a numeric 0.NNN.
$ws1.Range("H2:H6").NumberFormat = "@"

for ($i = 0; $i -lt $activeSignals.Length; $i++) {
    $row = $activeSignals[$i]
    $r = 2 + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws1.Cells.Item($r, 1 + $c).Value = $row[$c]
    }
}

# ---------------------------------------------------------------------
# Sheet 2: "Summary Dashboard"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary Dashboard")

$ws2.Range("B4").Value = 5       # Active Signals
$ws2.Range("B5").Value = 6       # BUY Signals
$ws2.Range("B6").Value = 9       # SELL Signals

# Average Confidence is stored as literal text (e.g. "77.6%"); force Text
# format so it isn't auto-converted into a numeric percentage.
$ws2.Range("B7").NumberFormat = "@"
$ws2.Range("B7").Value = "77.6%" # Average Confidence


# Average Risk/Reward is also literal text (e.g. "2.09"), not a number.
$ws2.Range("B8").NumberFormat = "@"
$ws2.Range("B8").Value = "2.09"  # Average Risk/Reward
$ws2.Range("B9").Value = "2025-07-28 16:19:05" # Last Update

# ---------------------------------------------------------------------
# Sheet 3: "Signal History"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Signal History")

$signalHistory = @(
    @("2025-07-28 16:38", "AUDUSD", "SELL", 0.65446, 0.65779, 0.64612, 0.05, 0.75, 2.5,  "Pending"),
    @("2025-07-28 16:29", "NZDUSD", "SELL", 0.59027, 0.59471, 0.5813,  0.05, 0.68, 2.02, "Active"),
    @("2025-07-28 16:31", "EURUSD", "BUY",  1.10518, 1.10121, 1.11047, 0.03, 0.79, 1.33, "Pending"),
    @("2025-07-28 16:03", "NZDUSD", "SELL", 0.59092, 0.59559, 0.58222, 0.01, 0.84, 1.86, "Active"),
    @("2025-07-28 15:53", "USDJPY", "SELL", 148.7591, 148.97191, 148.00594, 0.08, 0.66, 3.54, "Active"),
    @("2025-07-28 16:36", "NZDUSD", "BUY",  0.58971, 0.58547, 0.5941,  0.05, 0.91, 1.04, "Filled"),
    @("2025-07-28 16:34", "NZDUSD", "BUY",  0.59032, 0.5938,  0.58405, 0.02, 0.87, 1.8,  "Active"),
    @("2025-07-28 16:19", "AUDUSD", "SELL", 0.65982, 0.6634100000000001, 0.65043, 0.01, 0.92, 2.62, "Filled"),
    @("2025-07-28 16:34", "EURUSD", "BUY",  1.10384, 1.09903, 1.10812, 0.08, 0.75, 0.89, "Pending"),
    @("2025-07-28 15:54", "NZDUSD", "BUY",  0.59185, 0.58911, 0.60011, 0.03, 0.65, 3.02, "Filled"),
    @("2025-07-28 16:11", "USDJPY", "BUY",  150.20715, 149.90187, 150.6095, 0.07000000000000001, 0.67, 1.32, "Filled"),
    @("2025-07-28 16:06", "NZDUSD", "SELL", 0.58912, 0.59309, 0.5834,  0.09, 0.66, 1.44, "Filled"),
    @("2025-07-28 16:07", "USDCHF", "SELL", 0.88256, 0.88551, 0.87266, 0.09, 0.82, 3.36, "Active"),
    @("2025-07-28 16:29", "USDCAD", "SELL", 1.36602, 1.36135, 1.37252, 0.02, 0.75, 1.39, "Pending"),
    @("2025-07-28 16:44", "NZDUSD", "SELL", 0.58751, 0.5906400000000001, 0.57763, 0.05, 0.92, 3.16, "Filled")
)

for ($i = 0; $i -lt $signalHistory.Length; $i++) {
    $row = $signalHistory[$i]
    $r = 2 + $i
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws3.Cells.Item($r, 1 + $c).Value = $row[$c]
    }
}
